$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = "[name=`"Muelsyse`"]What do you think of the whole 'Diαbolic Crisis?'`n"
$ws.Range("C25").Value = "[name=`"Mayer`"]But the biggest impression I have of her is she's really cautious, going from what I know about the 'Diαbolic Crisis.'`n"
$ws.Range("C27").Value = "[name=`"Mayer`"]There were rumors among my coworkers that the 'Diαbolic Crisis' was some kinda conspiracy. I don't know that much about that stuff, though.`n"
$ws.Range("C32").Value = "[name=`"Muelsyse`"]But well, there's one thing I can tell you. Saria didn't leave Rhine Lab just because of the 'Diαbolic Crisis.' It's just one of a whole string of events that led to it, and there's more to it.`n"
$ws.Range("C42").Value = "[name=`"Muelsyse`"]But that whole string of events, the 'Diαbolic Crisis' included, meant that someone was trying to cross that line. They might even have crossed it already.`n"
$ws.Range("C48").Value = "[name=`"Muelsyse`"]You see? You're honest and kind, but you couldn't care less. You care more about your research, your 'Mibo.'`n"
$ws.Range("C123").Value = "[name=`"Muelsyse`"]Hehe, you think it's reckless, and you'd be right under normal circumstances. It's a risky move if he didn't have his 'prestige,' or perhaps 'influence' is the better word to use here.`n"
$ws.Range("C126").Value = "[name=`"Muelsyse`"]That's because what he said was also a threat. He was telling her, 'You could tell him, but it's not going to do you any good.'`n"
